# Heatmap and cell abundance profile
# Update the colour codes used for the Subantarctic zone ("SA-Sc") on the
# "cycle" sheet, and move the active sheet/selection focus from "Sheet3"
# to the "cycle" sheet.

$wb = $excel.ActiveWorkbook

$wsCycle = $wb.Worksheets.Item("cycle")

# Update the two colour values associated with the Subantarctic (SA-Sc)
# row on the "cycle" sheet - set B6 first so the new shared strings are
# emitted in the same order as the target workbook (#0D47A1 before
# #2196F3).
$wsCycle.Range("B6").Value = "#0D47A1"
$wsCycle.Range("B5").Value = "#2196F3"

# Move the active tab / selection from "Sheet3" to "cycle".
$wsCycle.Activate()
$wsCycle.Range("I11").Select()
